$wb = $excel.ActiveWorkbook

# --- Sheet "Metadatos": statistical operation updated ---
$wsMeta = $wb.Worksheets.Item("Metadatos")
$wsMeta.Range("F3").Value = "30324 Estimación de Defunciones Semanales"
$wsMeta.Range("F4").Value = "30324 Estimación de Defunciones Semanales"
$wsMeta.Range("F5").Value = "30324 Estimación de Defunciones Semanales"
$wsMeta.Range("F6").Value = "30324 Estimación de Defunciones Semanales"
$wsMeta.Range("F7").Value = "30324 Estimación de Defunciones Semanales"

# --- Sheet "Datos": strip accents from 6 CCAA names ---
$wsDatos = $wb.Worksheets.Item("Datos")
$wsDatos.Range("A3").Value  = "Andalucia"
$wsDatos.Range("A4").Value  = "Aragon"
$wsDatos.Range("A9").Value  = "Castilla y Leon"
$wsDatos.Range("A11").Value = "Cataluna"
$wsDatos.Range("A16").Value = "Murcia, Region de"
$wsDatos.Range("A18").Value = "Pais Vasco"

# --- Sheet "Clasificación": reclassify from Experimental to Sociedad ---
$wsClas = $wb.Worksheets.Item("Clasificación")
$wsClas.Range("B28").Value = 1
$wsClas.Range("B29").Value = 0

# --- Selection / view state updates ---
$wsDescripcion = $wb.Worksheets.Item("Descripción")
$wsDescripcion.Range("H19").Select()
$wsDatos.Range("A21").Select()
$wsMeta.Range("F6").Select()
$wsClas.Range("D27").Select()
